# Jan 5 - Updates
#
# Renames the two tabs ("Active Campaign" -> "Version 1", "Sheet1" ->
# "Active Campaign") and flips which of them is the active/selected tab,
# along with each sheet's current selection.

$wb = $excel.ActiveWorkbook

# Physical sheet order matches the <sheets> declaration order / r:id order:
#   Worksheets.Item(1) == rId1 (was "Active Campaign")
#   Worksheets.Item(2) == rId2 (was "Sheet1")
$sheetVersion1 = $wb.Worksheets.Item(1)
$sheetActiveCampaign = $wb.Worksheets.Item(2)

# Rename the first sheet out of the way before renaming the second sheet
# into "Active Campaign", otherwise the name would collide.
$sheetVersion1.Name = "Version 1"
$sheetActiveCampaign.Name = "Active Campaign"

# Update the selection on "Version 1" (no longer the active tab).
$sheetVersion1.Activate()
$sheetVersion1.Range("C27").Select()

# Update the selection on "Active Campaign" and leave it active/selected,
# matching the new tabSelected/activeTab state.
$sheetActiveCampaign.Activate()
$sheetActiveCampaign.Range("C78").Select()
